{"js": "// Update the worksheet date and the 25 two-digit-by-two-digit\n// multiplication problems to the newly generated set.\nconst replacements = [\n  [\"2024-09-26 Thursday\", \"2024-09-27 Friday\"],\n  [\"75\\u00d731=2325\", \"73\\u00d715=1095\"],\n  [\"79\\u00d779=6241\", \"23\\u00d711=253\"],\n  [\"26\\u00d743=1118\", \"28\\u00d739=1092\"],\n  [\"68\\u00d736=2448\", \"60\\u00d736=2160\"],\n  [\"51\\u00d739=1989\", \"67\\u00d739=2613\"],\n  [\"35\\u00d773=2555\", \"36\\u00d783=2988\"],\n  [\"64\\u00d757=3648\", \"93\\u00d791=8463\"],\n  [\"24\\u00d785=2040\", \"21\\u00d712=252\"],\n  [\"53\\u00d739=2067\", \"37\\u00d758=2146\"],\n  [\"93\\u00d747=4371\", \"45\\u00d733=1485\"],\n  [\"88\\u00d721=1848\", \"79\\u00d797=7663\"],\n  [\"46\\u00d749=2254\", \"24\\u00d744=1056\"],\n  [\"40\\u00d782=3280\", \"34\\u00d776=2584\"],\n  [\"51\\u00d735=1785\", \"63\\u00d771=4473\"],\n  [\"45\\u00d759=2655\", \"92\\u00d738=3496\"],\n  [\"72\\u00d776=5472\", \"16\\u00d717=272\"],\n  [\"61\\u00d775=4575\", \"89\\u00d786=7654\"],\n  [\"87\\u00d769=6003\", \"25\\u00d724=600\"],\n  [\"99\\u00d759=5841\", \"95\\u00d723=2185\"],\n  [\"91\\u00d725=2275\", \"14\\u00d722=308\"],\n  [\"37\\u00d773=2701\", \"25\\u00d761=1525\"],\n  [\"21\\u00d718=378\", \"51\\u00d732=1632\"],\n  [\"47\\u00d789=4183\", \"21\\u00d725=525\"],\n  [\"55\\u00d782=4510\", \"26\\u00d714=364\"],\n  [\"30\\u00d746=1380\", \"51\\u00d780=4080\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit-by-two-digit\n# multiplication problems to the newly generated set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-26 Thursday\", \"2024-09-27 Friday\"),\n    @(\"75\u00d731=2325\", \"73\u00d715=1095\"),\n    @(\"79\u00d779=6241\", \"23\u00d711=253\"),\n    @(\"26\u00d743=1118\", \"28\u00d739=1092\"),\n    @(\"68\u00d736=2448\", \"60\u00d736=2160\"),\n    @(\"51\u00d739=1989\", \"67\u00d739=2613\"),\n    @(\"35\u00d773=2555\", \"36\u00d783=2988\"),\n    @(\"64\u00d757=3648\", \"93\u00d791=8463\"),\n    @(\"24\u00d785=2040\", \"21\u00d712=252\"),\n    @(\"53\u00d739=2067\", \"37\u00d758=2146\"),\n    @(\"93\u00d747=4371\", \"45\u00d733=1485\"),\n    @(\"88\u00d721=1848\", \"79\u00d797=7663\"),\n    @(\"46\u00d749=2254\", \"24\u00d744=1056\"),\n    @(\"40\u00d782=3280\", \"34\u00d776=2584\"),\n    @(\"51\u00d735=1785\", \"63\u00d771=4473\"),\n    @(\"45\u00d759=2655\", \"92\u00d738=3496\"),\n    @(\"72\u00d776=5472\", \"16\u00d717=272\"),\n    @(\"61\u00d775=4575\", \"89\u00d786=7654\"),\n    @(\"87\u00d769=6003\", \"25\u00d724=600\"),\n    @(\"99\u00d759=5841\", \"95\u00d723=2185\"),\n    @(\"91\u00d725=2275\", \"14\u00d722=308\"),\n    @(\"37\u00d773=2701\", \"25\u00d761=1525\"),\n    @(\"21\u00d718=378\", \"51\u00d732=1632\"),\n    @(\"47\u00d789=4183\", \"21\u00d725=525\"),\n    @(\"55\u00d782=4510\", \"26\u00d714=364\"),\n    @(\"30\u00d746=1380\", \"51\u00d780=4080\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $new, \"wdReplaceAll\")\n}\n"}
